{"js": "// Update the worksheet date title and the 25 division problems/answers in\n// the 5x5 grid of filled table rows (rows 0, 4, 8, 12, 16 of the 20-row\n// table; the rows in between are blank \"work space\" rows).\n\nconst body = context.document.body;\n\n// --- Title paragraph: date line ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(\"2023-08-11 Friday\") !== -1) {\n    p.insertText(\"2023-08-12 Saturday\", \"Replace\");\n    break;\n  }\n}\n\n// --- Table cells ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Maps \"tableRow,col\" -> new text, for the 5 populated rows of the table.\nconst updates = {\n  \"0,0\": \"79\u00f79=8, 7\",\n  \"0,1\": \"50\u00f72=25, 0\",\n  \"0,2\": \"65\u00f75=13, 0\",\n  \"0,3\": \"15\u00f75=3, 0\",\n  \"0,4\": \"54\u00f73=18, 0\",\n\n  \"4,0\": \"72\u00f75=14, 2\",\n  \"4,1\": \"13\u00f77=1, 6\",\n  \"4,2\": \"42\u00f72=21, 0\",\n  \"4,3\": \"16\u00f75=3, 1\",\n  \"4,4\": \"18\u00f78=2, 2\",\n\n  \"8,0\": \"22\u00f76=3, 4\",\n  \"8,1\": \"71\u00f76=11, 5\",\n  \"8,2\": \"21\u00f76=3, 3\",\n  \"8,3\": \"89\u00f76=14, 5\",\n  \"8,4\": \"25\u00f76=4, 1\",\n\n  \"12,0\": \"12\u00f76=2, 0\",\n  \"12,1\": \"57\u00f79=6, 3\",\n  \"12,2\": \"95\u00f75=19, 0\",\n  \"12,3\": \"82\u00f74=20, 2\",\n  \"12,4\": \"46\u00f79=5, 1\",\n\n  \"16,0\": \"49\u00f75=9, 4\",\n  \"16,1\": \"32\u00f76=5, 2\",\n  \"16,2\": \"34\u00f72=17, 0\",\n  \"16,3\": \"49\u00f74=12, 1\",\n  \"16,4\": \"48\u00f77=6, 6\",\n};\n\nfor (const key in updates) {\n  const [rowStr, colStr] = key.split(\",\");\n  const row = parseInt(rowStr, 10);\n  const col = parseInt(colStr, 10);\n  const cell = table.getCell(row, col);\n  // Use `.value` (maps to the cell Range's Text property) rather than\n  // `body.insertText(...)`, so the existing run formatting (font/size) on\n  // the cell's text is preserved instead of being reset.\n  cell.value = updates[key];\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date title and the 25 division problems/answers in\n# the 5x5 grid of filled table rows (table rows 1, 5, 9, 13, 17 of the\n# 20-row table; the rows in between are blank \"work space\" rows).\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: date line ---\n$d.Paragraphs.Item(1).Range.Text = \"2023-08-12 Saturday\"\n\n# --- Table cells ---\n$t = $d.Tables.Item(1)\n\n# Row (1-based, only the 5 populated rows) -> 5 column values.\n$values = @{\n    1  = @(\"79\u00f79=8, 7\", \"50\u00f72=25, 0\", \"65\u00f75=13, 0\", \"15\u00f75=3, 0\", \"54\u00f73=18, 0\")\n    5  = @(\"72\u00f75=14, 2\", \"13\u00f77=1, 6\", \"42\u00f72=21, 0\", \"16\u00f75=3, 1\", \"18\u00f78=2, 2\")\n    9  = @(\"22\u00f76=3, 4\", \"71\u00f76=11, 5\", \"21\u00f76=3, 3\", \"89\u00f76=14, 5\", \"25\u00f76=4, 1\")\n    13 = @(\"12\u00f76=2, 0\", \"57\u00f79=6, 3\", \"95\u00f75=19, 0\", \"82\u00f74=20, 2\", \"46\u00f79=5, 1\")\n    17 = @(\"49\u00f75=9, 4\", \"32\u00f76=5, 2\", \"34\u00f72=17, 0\", \"49\u00f74=12, 1\", \"48\u00f77=6, 6\")\n}\n\nforeach ($row in $values.Keys) {\n    $cols = $values[$row]\n    for ($c = 0; $c -lt $cols.Length; $c++) {\n        $t.Cell($row, $c + 1).Range.Text = $cols[$c]\n    }\n}\n"}
